$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.77"
$ws.Range("E2").Value = "'1.38%"
$ws.Range("D3").Value = "'39.52"
$ws.Range("E3").Value = "'10.64%"
$ws.Range("D4").Value = "'5.105"
$ws.Range("E4").Value = "'1.32%"
$ws.Range("D5").Value = "'0.08150"
$ws.Range("E5").Value = "'3.15%"
$ws.Range("D6").Value = "'1.989"
$ws.Range("E6").Value = "'7.65%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.176"
$ws.Range("E7").Value = "'1.76%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.898"
$ws.Range("E8").Value = "'1.47%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9284"
$ws.Range("E9").Value = "'1.00%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1409"
$ws.Range("E10").Value = "'5.36%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1945"
$ws.Range("E11").Value = "'2.17%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09194"
$ws.Range("E12").Value = "'0.74%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03503"
$ws.Range("E13").Value = "'0.99%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09816"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001402"
$ws.Range("E15").Value = "'0.18%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006003"
$ws.Range("E16").Value = "'-2.25%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.947"
$ws.Range("E17").Value = "'5.92%"
$ws.Range("D18").Value = "'3.433"
$ws.Range("E18").Value = "'2.66%"
$ws.Range("D19").Value = "'0.3453"
$ws.Range("E19").Value = "'0.36%"
$ws.Range("D20").Value = "'0.1304"
$ws.Range("E20").Value = "'-0.47%"
$ws.Range("D21").Value = "'4.816"
$ws.Range("E21").Value = "'-6.62%"
$ws.Range("E22").Value = "'19.45%"
$ws.Range("D23").Value = "'0.04477"
$ws.Range("E23").Value = "'1.65%"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'0.65%"
$ws.Range("E25").Value = "'-9.66%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'0.07%"
$ws.Range("D39").Value = "'0.02114"
$ws.Range("E39").Value = "'9.12%"
$ws.Range("D40").Value = "'0.05160"
$ws.Range("D41").Value = "'0.007482"
$ws.Range("E41").Value = "'-1.71%"
$ws.Range("E42").Value = "'-0.22%"
$ws.Range("D43").Value = "'0.1366"
$ws.Range("E43").Value = "'1.79%"
$ws.Range("D44").Value = "'0.002132"
$ws.Range("E44").Value = "'-0.86%"
$ws.Range("D45").Value = "'0.009676"
$ws.Range("E45").Value = "'-4.90%"
$ws.Range("D46").Value = "'0.00006315"
$ws.Range("E46").Value = "'2.45%"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("E48").Value = "'1.94%"
$ws.Range("D49").Value = "'0.001602"
$ws.Range("E49").Value = "'-3.50%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.05%"
